# Junction_Flooding_234.xlsx edit:
#   - Row 5 values are replaced with a "custom accuracy" (2 decimal place) re-run
#     of the model, slightly different from a pure numeric rounding of the old
#     values.
#   - Row 6 (the extra simulation timestep) is removed entirely, shrinking the
#     used range from A1:AH6 down to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (B5:AH5) with the newly-computed, lower-precision values.
$ws.Range("B5").Value = 12.55
$ws.Range("C5").Value = 9.11
$ws.Range("D5").Value = 0.62
$ws.Range("E5").Value = 26.47
$ws.Range("F5").Value = 21.4
$ws.Range("G5").Value = 9.12
$ws.Range("H5").Value = 35.45
$ws.Range("I5").Value = 14.31
$ws.Range("J5").Value = 6.28
$ws.Range("K5").Value = 9.47
$ws.Range("L5").Value = 10.48
$ws.Range("M5").Value = 11.63
$ws.Range("N5").Value = 3.02
$ws.Range("O5").Value = 8.98
$ws.Range("P5").Value = 13.4
$ws.Range("Q5").Value = 7.75
$ws.Range("R5").Value = 0.4
$ws.Range("S5").Value = 0.25
$ws.Range("T5").Value = 135.53
$ws.Range("U5").Value = 26.18
$ws.Range("V5").Value = 8.86
$ws.Range("W5").Value = 17.87
$ws.Range("X5").Value = 9.5
$ws.Range("Y5").Value = 1.23
$ws.Range("Z5").Value = 17.01
$ws.Range("AA5").Value = 7.66
$ws.Range("AB5").Value = 6.44
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 11.14
$ws.Range("AE5").Value = 0.12
$ws.Range("AF5").Value = 31.35
$ws.Range("AG5").Value = 4.79
$ws.Range("AH5").Value = 10.67

# Remove row 6 completely (shifts dimension from A1:AH6 to A1:AH5).
$ws.Rows.Item(6).Delete()
